$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: status changes from "Ativo" to "Inativo"
$ws.Range("D3").Value = "Inativo"

# Row 4: validade_em_anos changes from 3 to 4
$ws.Range("C4").Value = 4

# New row 7: Primeiros Socorros | Engenheiro Quimico | 2 | Ativo
$ws.Range("A7").Value = "Primeiros Socorros"
$ws.Range("B7").Value = "Engenheiro Quimico"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "Ativo"

# Copy formatting for column A (name) from an existing data row
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Apply center/wrap formatting (matching the rest of the table) to the new cells
$rng = $ws.Range("B7:D7")
$rng.HorizontalAlignment = -4108
$rng.WrapText = $true

# Give the new row the same height used throughout the sheet
$ws.Rows.Item(7).RowHeight = 24.05

# Update selection to the last edited cell
$ws.Range("D7").Select() | Out-Null
